$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a numeric-looking string value need NumberFormat forced to
# Text ("@") first, otherwise Excel auto-converts the literal text into a
# floating point number (e.g. "1.170" -> 1.17), which would both lose the
# trailing zero and change the cell from Text to Number type.
$textCells = @('D5', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D24', 'D25', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D46', 'D47', 'D48', 'D49', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply every cell value change from the source diff.
$ws.Range('D2').Value = '30.199.60'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').Value = '2.080.12'
$ws.Range('E3').Value = '  -1.73%  '
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('D5').Value = '338.99'
$ws.Range('E6').Value = '  -0.45%  '
$ws.Range('D7').Value = '0.5268'
$ws.Range('E7').Value = '  +1.21%  '
$ws.Range('D8').Value = '0.4361'
$ws.Range('D9').Value = '54.85'
$ws.Range('E9').Value = '  +0.42%  '
$ws.Range('D10').Value = '0.09339'
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('D11').Value = '1.170'
$ws.Range('E11').Value = '  -0.84%  '
$ws.Range('D12').Value = '24.48'
$ws.Range('E12').Value = '  -2.84%  '
$ws.Range('D13').Value = '8.458'
$ws.Range('E13').Value = '  +0.12%  '
$ws.Range('D14').Value = '6.842'
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').Value = '2.035.32'
$ws.Range('E15').Value = '  -2.92%  '
$ws.Range('D16').Value = '100.30'
$ws.Range('E16').Value = '  -1.98%  '
$ws.Range('D17').Value = '0.00001159'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('D18').Value = '1.005'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').Value = '20.86'
$ws.Range('E19').Value = '  -3.18%  '
$ws.Range('D20').Value = '0.06678'
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('D21').Value = '6.300'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').Value = '1.003'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('D23').Value = '30.208.68'
$ws.Range('E23').Value = '  +0.79%  '
$ws.Range('D24').Value = '12.37'
$ws.Range('E24').Value = '  -3.03%  '
$ws.Range('D25').Value = '2.316'
$ws.Range('E25').Value = '  -0.65%  '
$ws.Range('E26').Value = '  -1.71%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '162.36'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '6.795'
$ws.Range('E28').Value = '  +4.19%  '
$ws.Range('D29').Value = '2.481'
$ws.Range('E29').Value = '  -2.73%  '
$ws.Range('D30').Value = '133.24'
$ws.Range('E30').Value = '  -0.57%  '
$ws.Range('D31').Value = '1.127'
$ws.Range('E31').Value = '  -2.63%  '
$ws.Range('D32').Value = '1.658'
$ws.Range('E32').Value = '  -7.22%  '
$ws.Range('D33').Value = '0.1046'
$ws.Range('E33').Value = '  -1.09%  '
$ws.Range('D34').Value = '6.231'
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('D35').Value = '3.906'
$ws.Range('E35').Value = '  -1.52%  '
$ws.Range('D36').Value = '0.02603'
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').Value = '9.849'
$ws.Range('E37').Value = '  -8.95%  '
$ws.Range('D38').Value = '0.06709'
$ws.Range('E38').Value = '  -2.30%  '
$ws.Range('D39').Value = '0.6944'
$ws.Range('E39').Value = '  -1.21%  '
$ws.Range('D40').Value = '12.51'
$ws.Range('E40').Value = '  -1.35%  '
$ws.Range('D41').Value = '1.331'
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').Value = '0.2199'
$ws.Range('E42').Value = '  -2.21%  '
$ws.Range('D43').Value = '0.6712'
$ws.Range('E43').Value = '  -1.96%  '
$ws.Range('D44').Value = '2.360'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('E45').Value = '  -1.48%  '
$ws.Range('D46').Value = '1.002'
$ws.Range('E46').Value = '  -0.37%  '
$ws.Range('D47').Value = '1.314'
$ws.Range('E47').Value = '  +5.43%  '
$ws.Range('D48').Value = '3.625'
$ws.Range('E48').Value = '  -0.35%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.00000000352'
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').Value = '1.210'
$ws.Range('E50').Value = '  +2.14%  '
$ws.Range('E51').Value = '  -1.44%  '
